$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$pr = $last.Range
Write-Output ("last para range: " + $pr.Start + "-" + $pr.End)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00777F1F" w:rsidRPr="004C6456" w:rsidRDefault="00777F1F" w:rsidP="004C6456"/><w:sectPr w:rsidR="00777F1F" w:rsidRPr="004C6456" w:rsidSect="00140843"><w:footerReference w:type="even" r:id="rId8"/><w:footerReference w:type="default" r:id="rId9"/><w:headerReference w:type="first" r:id="rId10"/><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1008" w:right="1800" w:bottom="1440" w:left="1800" w:header="720" w:footer="864" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="326"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$res = $pr.InsertXML($frag)
Write-Output ("result: " + $res)
